# Generate Report for Handback
#
# Refreshes the localization-status report after a successful handback:
#   - "Status" cells that said "Ready for handoff" now read
#     "Handed back: in sync with en-US" (Overview!E2/F2 and the per-locale
#     "Status" column on the zh-cn / de-de sheets all shared that string).
#   - The "Latest Handback DateTime" for zh-cn / de-de is refreshed to the
#     handback timestamp.
#   - The stale "Error Detail" message (handback file was behind the
#     latest) is cleared now that the handback is in sync.
#   - The "Status" / "Error Detail" columns are widened / narrowed to fit
#     the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status summary columns (E, F)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# zh-cn sheet: Status (C2), Latest Handback DateTime (K2), Error Detail (P2)
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-18 20:50:54"
$wsZhCn.Range("P2").Value = ""

# de-de sheet: Status (C2), Latest Handback DateTime (K2), Error Detail (P2)
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-18 20:51:05"
$wsDeDe.Range("P2").Value = ""

# Column width adjustments (character units) to fit the new text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
